$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update H8: expenses for 2018 (now actual, was placeholder 0) ---
$ws.Range("H8").Value = -1330.14

# --- Update E8: escrow payment formula now explicit for 2018 (actual), like rows 2,6,7 ---
$ws.Range("E8").Formula = "=-2234.42-C8-D8"

# --- Highlight row 8 (2018) as "actual" data, matching rows 2-7 ---
$ws.Range("A8:P8").Interior.Color = 5296274

# --- Update selection / view to reflect new active cell ---
$ws.Range("E9").Select()
